$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price column so values like "12.00" or "0.0000246"
# are preserved verbatim instead of being re-interpreted as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '68.667.83'
$ws.Range('E2').Value = '  -0.94%  '
$ws.Range('D3').Value = '3.914.40'
$ws.Range('E3').Value = '  +3.34%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '601.86'
$ws.Range('E5').Value = '  -0.33%  '
$ws.Range('D6').Value = '165.38'
$ws.Range('E6').Value = '  +0.12%  '
$ws.Range('D7').Value = '3.914.67'
$ws.Range('E7').Value = '  +3.44%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('E9').Value = '  -2.15%  '
$ws.Range('E10').Value = '  -3.98%  '
$ws.Range('D11').Value = '6.35'
$ws.Range('E11').Value = '  +0.38%  '
$ws.Range('E12').Value = '  -0.47%  '
$ws.Range('B13').Value = 'ShibaInu'
$ws.Range('C13').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D13').Value = '0.0000246'
$ws.Range('E13').Value = '  -0.54%  '
$ws.Range('B14').Value = 'Avalanche'
$ws.Range('C14').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D14').Value = '36.98'
$ws.Range('E14').Value = '  -1.72%  '
$ws.Range('D15').Value = '4.565.95'
$ws.Range('D16').Value = '3.944.13'
$ws.Range('E16').Value = '  +4.26%  '
$ws.Range('D17').Value = '68.811.60'
$ws.Range('E17').Value = '  -0.86%  '
$ws.Range('E18').Value = '  -0.34%  '
$ws.Range('E19').Value = '  -1.55%  '
$ws.Range('D20').Value = '17.07'
$ws.Range('E20').Value = '  -3.18%  '
$ws.Range('D21').Value = '11.13'
$ws.Range('E21').Value = '  -1.77%  '
$ws.Range('D22').Value = '484.03'
$ws.Range('E22').Value = '  -1.83%  '
$ws.Range('E23').Value = '  -1.08%  '
$ws.Range('E24').Value = '  +11.67%  '
$ws.Range('D25').Value = '84.49'
$ws.Range('E25').Value = '  -0.48%  '
$ws.Range('E26').Value = '  -1.41%  '
$ws.Range('D27').Value = '12.00'
$ws.Range('E27').Value = '  -2.57%  '
$ws.Range('D28').Value = '10.08'
$ws.Range('E28').Value = '  -0.43%  '
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('E30').Value = '  -1.92%  '
$ws.Range('D31').Value = '4.063.26'
$ws.Range('E31').Value = '  +3.33%  '
$ws.Range('E32').Value = '  -3.65%  '
$ws.Range('E33').Value = '  -2.18%  '
$ws.Range('D34').Value = '31.99'
$ws.Range('E34').Value = '  +0.23%  '
$ws.Range('D35').Value = '3.856.59'
$ws.Range('E35').Value = '  +3.26%  '
$ws.Range('E36').Value = '  -1.41%  '
$ws.Range('E37').Value = '  +1.99%  '
$ws.Range('E38').Value = '  -0.45%  '
$ws.Range('D39').Value = '5.88'
$ws.Range('E39').Value = '  -1.42%  '
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').Value = '  +0.10%  '
$ws.Range('D41').Value = '3.08'
$ws.Range('E41').Value = '  +0.57%  '
$ws.Range('E42').Value = '  -2.58%  '
$ws.Range('D43').Value = '432.28'
$ws.Range('E43').Value = '  +1.97%  '
$ws.Range('D44').Value = '48.50'
$ws.Range('E44').Value = '  -0.04%  '
$ws.Range('E45').Value = '  -0.86%  '
$ws.Range('E46').Value = '  +0.00%  '
$ws.Range('D47').Value = '8.42'
$ws.Range('E47').Value = '  -0.20%  '
$ws.Range('D48').Value = '26.45'
$ws.Range('E48').Value = '  +9.71%  '
$ws.Range('D49').Value = '2.821.15'
$ws.Range('E49').Value = '  +0.24%  '
$ws.Range('D50').Value = '141.86'
$ws.Range('E50').Value = '  +0.13%  '
$ws.Range('B51').Value = 'FLOKI'
$ws.Range('C51').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D51').Value = '0.000265'
$ws.Range('E51').Value = '  +17.47%  '
